$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stocks")
$ws.Select()
